$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at Excel row 403, shifting the existing rows 403-420
# (the remainder of the "Espinaca" price series for Terminal La Palmera de
# La Serena) down to 404-421.
$ws.Rows.Item(403).Insert()

# Populate the newly inserted row with the new weekly price record.
$ws.Cells.Item(403, 1).Value  = 8
$ws.Cells.Item(403, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(403, 3).Value  = "Coquimbo"
$ws.Cells.Item(403, 4).Value  = 45075
$ws.Cells.Item(403, 5).Value  = 4
$ws.Cells.Item(403, 6).Value  = 100112012
$ws.Cells.Item(403, 7).Value  = "Espinaca"
$ws.Cells.Item(403, 8).Value  = "Sin especificar"
$ws.Cells.Item(403, 9).Value  = "Primera"
$ws.Cells.Item(403, 10).Value = 700
$ws.Cells.Item(403, 11).Value = 400
$ws.Cells.Item(403, 12).Value = 500
$ws.Cells.Item(403, 13).Value = 450
$ws.Cells.Item(403, 14).Value = "$/atado 300 a 500 gramos"
$ws.Cells.Item(403, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(403, 16).Value = 900
$ws.Cells.Item(403, 17).Value = 0.5
$ws.Cells.Item(403, 18).Value = "Hortaliza"
